$wb = $excel.ActiveWorkbook

# --- "Daily Orders": a brand-new order came in, so it goes on top of the log ---
$ws = $wb.Worksheets.Item("Daily Orders")

# Push the existing data rows (2-10) down one row, leaving row 2 empty.
$ws.Range("A2:L2").Insert()

# Fill in the new order's details.
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "2026-01-13 19:11"
$ws.Range("C2").Value = "Sagar Borse"

# Phone number must stay text (it would otherwise be read as a plain
# number) - write it as a formula returning the text, then collapse the
# formula down to its literal value so the cell stores static text.
$ws.Range("D2").Formula = "=""7588930329"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

$ws.Range("E2").Value = "Test,"
$ws.Range("F2").Value = "Stainless Steel Grater x1"
$ws.Range("G2").Value = 40
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"

$excel.CutCopyMode = $false

# --- "Summary": refresh the daily roll-up now that order #10 exists ---
$sum = $wb.Worksheets.Item("Summary")
$sum.Range("A2").Value = 10
$sum.Range("B2").Value = 7
$sum.Range("G2").Value = 600
